# Apply the benchmark-table edit: update the single-column results table
# in place, row by row, using the Table/Cell object model so formatting
# (Times New Roman, sz 22) on each run is preserved automatically.

$d = $word.ActiveDocument
$t = $d.Tables(1)

# Map of 1-based row index -> new cell text.
$changes = [ordered]@{
    1  = "0M"
    2  = "0M"
    3  = "0M"
    4  = "7584"
    5  = "0.00001"
    6  = "0.00112"
    7  = "0.00019"
    8  = "0.00006"
    9  = "0.00035"
    10 = "0.00051"
    11 = "0.00064"
    12 = "1.72374"
    44 = "99.88"
    45 = "1.72"
    46 = "1476"
}

foreach ($rowIndex in $changes.Keys) {
    $cell = $t.Cell($rowIndex, 1)
    $cell.Range.Text = $changes[$rowIndex]
}
